$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(15.15329931950635,10.32029484289569,0,15.73922050295,36.97416411342906,3.644082794583784,0,23.82156964551764,7.745291093423609,0,12.73455317073982,16.60815068032149,18.45601016960284,23.51964366730065),
    @(14.74324394752008,10.15349901774341,0,15.77187840597638,37.01116878736758,3.646026061369065,0,23.91798183045355,7.737483507866232,0,12.7331782029122,16.52163129511129,18.50707063950911,23.57895703228353),
    @(14.48756323597926,10.04856579192023,0,15.79342641532775,37.0423792531146,3.647283316922919,0,23.98177142547627,7.732691224403032,0,12.73384558831676,16.470440958795,18.54017735440427,23.62075773690207),
    @(14.38254393331393,10.00520058433682,0,15.80258421615792,37.05723081509729,3.647811822928096,0,24.00891981598357,7.730739019980808,0,12.7344990291335,16.45008169946048,18.55411099610703,23.63914209007698),
    @(14.36506028171986,9.997964225599155,0,15.804127639629,37.05982566986776,3.647900558702433,0,24.01349742868504,7.730414907895335,0,12.73463061051839,16.44673176325069,18.55645141486623,23.64227625177362),
    @(14.48615004794924,10.04798335954695,0,15.79354839406432,37.04257091274938,3.647290379024887,0,23.9821328887034,7.732664893125493,0,12.73385285453862,16.47016433772322,18.540363475855,23.62100021204709),
    @(15.0128126754146,10.26332406338186,0,15.75017085973543,36.98516056775317,3.644739563066059,0,23.85385880150699,7.742598875711508,0,12.73376619346971,16.57792690372307,18.4732520803215,23.53897605592422),
    @(16.0082914788186,10.66437116370266,0,15.67694810819113,36.93998159194252,3.64024358996479,0,23.63879995954914,7.762080139944838,0,12.74552722920148,16.80389346340912,18.35553042478596,23.42096201424052),
    @(16.70916251313636,10.9445006758408,0,15.63032984738147,36.94788787230778,3.637245766331567,0,23.50310741553927,7.776378258541652,0,12.76135238669203,16.97787739135116,18.27744211311814,23.36053522278583),
    @(17.01994467765008,11.06847794358843,0,15.61067231800944,36.96039163581578,3.635947604785566,0,23.44623709266865,7.782876364927688,0,12.7700904570526,17.05854269021897,18.24372914874709,23.33878100274889),
    @(17.13636648978818,11.11490397609165,0,15.60345066669055,36.96640410291731,3.63546540109171,0,23.4254015660426,7.785335830246812,0,12.77361864126847,17.08928831154194,18.23122223538535,23.33136969769546),
    @(17.11135081993007,11.10492886935116,0,15.60499610384611,36.96505245351509,3.635568835750589,0,23.42985770960555,7.784806200364472,0,12.77284906574632,17.08265812837352,18.23390429604913,23.33292906620673),
    @(17.02954869560689,11.07230802842912,0,15.6100737378575,36.96086070319434,3.635907745832113,0,23.4445089024069,7.7830787331538,0,12.7703763406666,17.06106826422441,18.24269500207019,23.33815469492333),
    @(16.97927481970033,11.05225820329217,0,15.61321285922907,36.95845939704054,3.636116558484828,0,23.45357439380975,7.782020439767182,0,12.76889021706208,17.04786925033366,18.24811332879966,23.34146323680962),
    @(16.68868081954466,10.93632680797353,0,15.63164564160291,36.94724969988291,3.637331919601478,0,23.50692189576791,7.775953438877105,0,12.76081213117715,16.97263467350952,18.27968168423202,23.36207247571568),
    @(16.50827002485705,10.86430376542177,0,15.64334997893547,36.94265224134446,3.638094263780526,0,23.5408939724018,7.772229708648362,0,12.75624926758804,16.92685634301644,18.29951081084776,23.37618575844517),
    @(16.40375080678813,10.82255419676404,0,15.65022787023228,36.94084654854121,3.638538917859439,0,23.56089086888128,7.77008735636248,0,12.75376982964793,16.90067000582562,18.31108639307374,23.38484298130445),
    @(16.36823672707131,10.80836366703327,0,15.65258167596494,36.94037929308994,3.638690531868494,0,23.56773992538151,7.76936190653581,0,12.75295529931812,16.89182909438525,18.31503497597606,23.38786678570155),
    @(16.52755366985562,10.87200443818715,0,15.64208893977387,36.94305487267223,3.638012472303156,0,23.53723027154547,7.772626166595724,0,12.75671999973434,16.93171474259802,18.29738233722345,23.37462750975557),
    @(17.05361109143985,11.08190390213986,0,15.60857628695386,36.96205728029852,3.635807945469041,0,23.44018648037322,7.783586168151636,0,12.77109670611196,17.06740447139643,18.24010592331569,23.33659735484674),
    @(17.39000578872823,11.21603536729194,0,15.58796888692687,36.98192083277581,3.634421823597945,0,23.3808438239922,7.790741853522455,0,12.78176973121104,17.15723893667175,18.20418442205031,23.3165608661209),
    @(17.21117694916282,11.14473373442824,0,15.59884913307674,36.97063946370758,3.635156636292902,0,23.41214216146987,7.786923511194735,0,12.77595719989278,17.10919341180493,18.22321831680659,23.3268132685746),
    @(16.51883801328132,10.86852402918976,0,15.64265859170006,36.94287023422487,3.638049430377093,0,23.53888517956609,7.772446932693858,0,12.75650673361742,16.92951784679613,18.29834407324153,23.37533030175178),
    @(15.74382320679674,10.55831805655532,0,15.69549331952074,36.94498312962038,3.641406013775356,0,23.69306599862123,7.756811414691625,0,12.74107700486536,16.80389346340912,18.38589740330444,23.42096201424052)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $rowVals[$j]
    }
}